$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 385; existing rows 385..438 shift down to 386..439
$ws.Rows.Item(385).Insert()

# Populate the newly inserted row 385 with the new record's data
$ws.Cells.Item(385,1).Value = 4
$ws.Cells.Item(385,2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(385,3).Value = "Los Lagos"
$ws.Cells.Item(385,4).Value = 45127
$ws.Cells.Item(385,5).Value = 10
$ws.Cells.Item(385,6).Value = 100112021
$ws.Cells.Item(385,7).Value = "Ají"
$ws.Cells.Item(385,8).Value = "Inferno"
$ws.Cells.Item(385,9).Value = "Primera"
$ws.Cells.Item(385,10).Value = 75
$ws.Cells.Item(385,11).Value = 22000
$ws.Cells.Item(385,12).Value = 22000
$ws.Cells.Item(385,13).Value = 22000
$ws.Cells.Item(385,14).Value = "$/caja 10 kilos"
$ws.Cells.Item(385,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(385,16).Value = 2200
$ws.Cells.Item(385,17).Value = 10
$ws.Cells.Item(385,18).Value = "Hortaliza"
